$wb = $excel.ActiveWorkbook

# 展览 (Exhibition) sheet
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 1263
$ws1.Range("F5").Value = 5589
$ws1.Range("F6").Value = 1790
$ws1.Range("F7").Value = 1790
$ws1.Range("F8").Value = 6361
$ws1.Range("F10").Value = 1927
$ws1.Range("F17").Value = 50
$ws1.Range("F18").Value = 7935
$ws1.Range("F19").Value = 7935
$ws1.Range("F31").Value = 1758
$ws1.Range("F32").Value = 802
$ws1.Range("F33").Value = 375

# 演出 (Performance) sheet
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F8").Value = 442
$ws2.Range("F21").Value = 63

# 本地生活 (Local Life) sheet
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F4").Value = 686
$ws3.Range("F5").Value = 268

# 全部类型 (All Types) sheet
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 686
$ws4.Range("F5").Value = 1263
$ws4.Range("F10").Value = 5589
$ws4.Range("F11").Value = 268
$ws4.Range("F12").Value = 1790
$ws4.Range("F13").Value = 1790
$ws4.Range("F14").Value = 6361
$ws4.Range("F16").Value = 1927
$ws4.Range("F23").Value = 50
$ws4.Range("F24").Value = 7935
$ws4.Range("F25").Value = 7935
$ws4.Range("F36").Value = 1758
$ws4.Range("F37").Value = 802
$ws4.Range("F39").Value = 375
$ws4.Range("F49").Value = 63
